# Automatische test-sync: 2025-06-18 10:00:10
#
# Adds a new incoming-mail row (row 5) to the "Logs" sheet and the
# corresponding aggregated count row (row 5) to the "Dashboard" sheet,
# then extends the conditional formatting ranges and the chart series
# ranges so they keep covering all the data rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Logs sheet: append the new mail entry as row 5
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A5").Value = "Klacht over levering"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$logs.Range("D5").Value = "Klacht"
$logs.Range("F5").Value = "2025-06-18 09:30:12"
$logs.Range("G5").Value = "Nee"

# Extend the conditional formatting that highlights the "Categorie" (D)
# and "Beantwoord" (G) columns so it also covers the new row.
$catFormats = $logs.Range("D2:D4").FormatConditions
for ($i = 1; $i -le $catFormats.Count; $i++) {
    $catFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D5"))
}

$answeredFormats = $logs.Range("G2:G4").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G5"))
}

# ---------------------------------------------------------------
# 2. Dashboard sheet: add the aggregated count for the new category
# ---------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A5").Value = "Klacht"
$dashboard.Range("B5").Value = 1

# ---------------------------------------------------------------
# 3. Chart: extend the category/value series ranges to include row 5
# ---------------------------------------------------------------
$chartObj = $dashboard.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$5,Dashboard!`$B`$2:`$B`$5,1)"
